$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "Energy Supplied till 2nd oscillation"
$ws.Range("B4").Value = 6160
$ws.Range("C4").Value = 3916
$ws.Range("D4").Value = 2779
$ws.Range("E4").Value = 1861
